$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Interpolation_instructions")
$ws2 = $wb.Worksheets.Item("Trend_instructions")

# --- Trend_instructions (sheet2) data edits ---
# Order matters for shared-string table layout: hard_coal, brown_coal, then linear (on sheet1) last.
$ws2.Range("B3").Value = "hard_coal"
$ws2.Range("C3").Value = "1A1a_Electricity-public"
$ws2.Range("D3").Value = 1932

$ws2.Range("B4").Value = "brown_coal"
$ws2.Range("C4").Value = "1A1a_Electricity-public"

$ws2.Range("C5").Value = "1A1a_Electricity-autoproducer"

$ws2.Range("E2").Value = 1934

# --- Interpolation_instructions (sheet1) data edit ---
$ws1.Range("H2").Value = "linear"

# --- Column width for Trend_instructions column C (CEDS_sector) ---
$ws2.Columns.Item(3).ColumnWidth = 25.1666666666667

# --- Sheet view / selection / active tab updates ---
$ws1.Activate()
$ws1.Range("H2").Select()

$ws2.Activate()
$ws2.Range("H12").Select()
